$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("M2").Value = 8.533046666666666
$ws.Range("N2").Value = 25.59914
$ws.Range("O2").Value = 0.2932132236642383
$ws.Range("P2").Value = 0.2932132236642383
$ws.Range("Q2").Value = 1.174940794673333
$ws.Range("R2").Value = 10.57446715206
$ws.Range("S2").Value = 0.2932132236642383
$ws.Range("T2").Value = 0.2932132236642383

# Row 3
$ws.Range("O3").Value = 0.3119288965200195
$ws.Range("P3").Value = 0.3119288965200194
$ws.Range("S3").Value = 0.3119288965200195
$ws.Range("T3").Value = 0.3119288965200194

# Row 4
$ws.Range("O4").Value = 0.3948578798157423
$ws.Range("P4").Value = 0.3948578798157423
$ws.Range("S4").Value = 0.3948578798157423
$ws.Range("T4").Value = 0.3948578798157423
